$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Ciserano Italy" (sheet index 2)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ciserano Italy")

$ws.Range("E2").Value = 0.0361
$ws.Range("E3").Value = 0.0361
$ws.Range("E4").Value = 0.0361

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.012
$ws.Range("O4").Value = 0.00601666666666667
$ws.Range("P4").Value = 0.00601666666666667
$ws.Range("Q4").Value = 0.00601666666666667
$ws.Range("R4").Value = 0.01805
$ws.Range("S4").Value = 0.00601666666666667
$ws.Range("T4").Value = 0.00601666666666667
$ws.Range("U4").Value = 0.00601666666666667
$ws.Range("V4").Value = 0.01805
$ws.Range("W4").Value = 0.0722

$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 0.8571

# ---------------------------------------------------------------------------
# Sheet "Molndal Sweden" (sheet index 10)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Molndal Sweden")
$ws.Range("M5").Value = $null

# ---------------------------------------------------------------------------
# Sheet "Piedras Negras Fasco Mexico" (sheet index 12)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Piedras Negras Fasco Mexico")
$ws.Range("M7").Value = $null

# ---------------------------------------------------------------------------
# Sheet "Waldenburg Germany" (sheet index 16)
# Insert a new "Internal Fill Rate" / "Commit/Forecast" row (row 7),
# shifting the old "Manufacturing Voluntary Turnover" rows down by one,
# and update the E / O-W values on those shifted rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Waldenburg Germany")

$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

$ws.Rows(7).Insert()

$ws.Range("A7").Value = "PES"
$ws.Range("B7").Value = "PES EMEA"
$ws.Range("C7").Value = "Waldenburg Germany"
$ws.Range("D7").Value = "Internal Fill Rate"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "Commit/Forecast"
$ws.Range("G7").Value = $null
$ws.Range("H7").Value = $null
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = $null
$ws.Range("K7").Value = $null
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0

# Row 8 (previously row 7 pre-insert): Manufacturing Voluntary Turnover / PY Actual
$ws.Range("E8").Value = 0.0302
$ws.Range("O8").Value = 0.0093
$ws.Range("P8").Value = 0.0093
$ws.Range("Q8").Value = 0.0089
$ws.Range("R8").Value = 0.0275
$ws.Range("S8").Value = 0.0092
$ws.Range("V8").Value = 0.0092
$ws.Range("W8").Value = 0.037

# Row 9 (previously row 8 pre-insert): Manufacturing Voluntary Turnover / AOP
$ws.Range("E9").Value = 0.0302
$ws.Range("O9").Value = 0.00837
$ws.Range("P9").Value = 0.00837
$ws.Range("Q9").Value = 0.00801
$ws.Range("R9").Value = 0.02475
$ws.Range("S9").Value = 0.00828
$ws.Range("V9").Value = 0.00828
$ws.Range("W9").Value = 0.0333

# Row 10 (previously row 9 pre-insert): Manufacturing Voluntary Turnover / Commit-Forecast
$ws.Range("E10").Value = 0.0302
$ws.Range("G10").Value = 0.03
$ws.Range("J10").Value = 0.0302
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0.00503333333333333
$ws.Range("P10").Value = 0.00503333333333333
$ws.Range("Q10").Value = 0.00503333333333333
$ws.Range("R10").Value = 0.0151
$ws.Range("S10").Value = 0.00503333333333333
$ws.Range("T10").Value = 0.00503333333333333
$ws.Range("U10").Value = 0.00503333333333333
$ws.Range("V10").Value = 0.0151
$ws.Range("W10").Value = 0.0604

# ---------------------------------------------------------------------------
# Sheet "Changzhou Epc China" (sheet index 21)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Changzhou Epc China")
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = $null

# ---------------------------------------------------------------------------
# Sheet "Bangalore India" (sheet index 24)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bangalore India")
$ws.Range("E2").Value = 0.3774
$ws.Range("E3").Value = 0.3774
$ws.Range("E4").Value = 0.3774

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.4255
$ws.Range("O4").Value = 0.0629
$ws.Range("P4").Value = 0.0629
$ws.Range("Q4").Value = 0.0629
$ws.Range("R4").Value = 0.1887
$ws.Range("S4").Value = 0.0629
$ws.Range("T4").Value = 0.0629
$ws.Range("U4").Value = 0.0629
$ws.Range("V4").Value = 0.1887
$ws.Range("W4").Value = 0.7548

# ---------------------------------------------------------------------------
# Sheet "Black River Falls Wisconsin" (sheet index 25)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Black River Falls Wisconsin")
$ws.Range("E2").Value = 0.2632
$ws.Range("E3").Value = 0.2632
$ws.Range("E4").Value = 0.2632

$ws.Range("I4").Value = 0.25
$ws.Range("J4").Value = 0.2128
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0438666666666667
$ws.Range("P4").Value = 0.0438666666666667
$ws.Range("Q4").Value = 0.0438666666666667
$ws.Range("R4").Value = 0.1316
$ws.Range("S4").Value = 0.0438666666666667
$ws.Range("T4").Value = 0.0438666666666667
$ws.Range("U4").Value = 0.0438666666666667
$ws.Range("V4").Value = 0.1316
$ws.Range("W4").Value = 0.5264

# ---------------------------------------------------------------------------
# Sheet "Juarez FCDM" (sheet index 26)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Juarez FCDM")
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null

# ---------------------------------------------------------------------------
# Sheet "Mumbai India" (sheet index 27)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mumbai India")
$ws.Range("E2").Value = 0.8571
$ws.Range("E3").Value = 0.8571
$ws.Range("E4").Value = 0.8571

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 1.2
$ws.Range("O4").Value = 0.14285
$ws.Range("P4").Value = 0.14285
$ws.Range("Q4").Value = 0.14285
$ws.Range("R4").Value = 0.42855
$ws.Range("S4").Value = 0.14285
$ws.Range("T4").Value = 0.14285
$ws.Range("U4").Value = 0.14285
$ws.Range("V4").Value = 0.42855
$ws.Range("W4").Value = 1.7142

# ---------------------------------------------------------------------------
# Sheet "Noida India" (sheet index 28)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Noida India")
$ws.Range("E2").Value = 0.1182
$ws.Range("E3").Value = 0.1182
$ws.Range("E4").Value = 0.1182

$ws.Range("G4").Value = 0.0227
$ws.Range("I4").Value = 0.0233
$ws.Range("J4").Value = 0.0458
$ws.Range("K4").Value = 0.0233
$ws.Range("L4").Value = 0.05
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0732
$ws.Range("O4").Value = 0.0197
$ws.Range("P4").Value = 0.0197
$ws.Range("Q4").Value = 0.0197
$ws.Range("R4").Value = 0.0591
$ws.Range("S4").Value = 0.0197
$ws.Range("T4").Value = 0.0197
$ws.Range("U4").Value = 0.0197
$ws.Range("V4").Value = 0.0591
$ws.Range("W4").Value = 0.2364

# ---------------------------------------------------------------------------
# Sheet "Juarez Casa I" (sheet index 29)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Juarez Casa I")
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null
$ws.Range("O4").Value = 0.222216666666667
$ws.Range("P4").Value = 0.222216666666667
$ws.Range("Q4").Value = 0.222216666666667
$ws.Range("R4").Value = 0.66665
$ws.Range("S4").Value = 0.222216666666667
$ws.Range("T4").Value = 0.222216666666667
$ws.Range("U4").Value = 0.222216666666667
$ws.Range("V4").Value = 0.66665
$ws.Range("W4").Value = 2.6666
